$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Other")

# --- Insert a new row for Elevation data, right after the Open Street Map row (22) ---
$ws.Rows.Item(23).Insert()
$ws.Range("A23").Value = "elevation"
$ws.Range("B23").Value = "NASA Shuttle Radar Topography Mission (SRTMGL1)"
$ws.Range("F23").Value = "1 arc second"
$ws.Range("E22").Value = "once"
$ws.Range("E23").Value = "once"
$ws.Range("G23").Value = "elevation data"

# --- Rename the NDVI/EVI rows (still rows 18/19) to the more specific 250m variants,
#     then insert a new row for the 1km NDVI product right after them (becomes row 20) ---
$ws.Range("A18").Value = "ndvi_250m"
$ws.Range("A19").Value = "evi_250m"
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = "ndvi_1km"
$ws.Range("B20").Value = "MOD13A3 v006"
$ws.Range("E20").Value = "30-day"
$ws.Range("F20").Value = "1km"
$ws.Range("G20").Value = "slightly lower resolution NDVI product; sufficient for purposes of NO2 model"
$ws.Range("C20").Value = 20050101
$ws.Range("D20").Value = 20191231

# --- Update selection to match the saved workbook state ---
$ws.Activate()
$ws.Range("B21").Select()
